$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 112.466131
$ws.Range("H2").Value = 337.3983929999999
$ws.Range("I2").Value = 0.2948491373870378
$ws.Range("J2").Value = 0.2948491373870378
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.307106666666667
$ws.Range("N2").Value = 3.92132
$ws.Range("O2").Value = 0.01256263154946851
$ws.Range("P2").Value = 0.01256263154946851
$ws.Range("Q2").Value = 147.0052296043066
$ws.Range("R2").Value = 1323.04706643876
$ws.Range("S2").Value = 0.003704081075671977
$ws.Range("T2").Value = 0.003704081075671977

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 112.466131
$ws.Range("H3").Value = 337.3983929999999
$ws.Range("I3").Value = 0.2948491373870378
$ws.Range("J3").Value = 0.2948491373870378
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 80.22623699999998
$ws.Range("N3").Value = 240.678711
$ws.Range("O3").Value = 0.77105616682495
$ws.Range("P3").Value = 0.77105616682495
$ws.Range("Q3").Value = 9022.734480079043
$ws.Range("R3").Value = 81204.61032071139
$ws.Range("S3").Value = 0.2273452456652925
$ws.Range("T3").Value = 0.2273452456652925

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 112.466131
$ws.Range("H4").Value = 337.3983929999999
$ws.Range("I4").Value = 0.2948491373870378
$ws.Range("J4").Value = 0.2948491373870378
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 22.51385866666667
$ws.Range("N4").Value = 67.541576
$ws.Range("O4").Value = 0.2163812016255815
$ws.Range("P4").Value = 0.2163812016255815
$ws.Range("Q4").Value = 2532.046578120818
$ws.Range("R4").Value = 22788.41920308737
$ws.Range("S4").Value = 0.0637998106460734
$ws.Range("T4").Value = 0.0637998106460734

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 82.26089466666666
$ws.Range("H5").Value = 246.782684
$ws.Range("I5").Value = 0.2156609604819841
$ws.Range("J5").Value = 0.2156609604819841
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.307106666666667
$ws.Range("N5").Value = 3.92132
$ws.Range("O5").Value = 0.01256263154946851
$ws.Range("P5").Value = 0.01256263154946851
$ws.Range("Q5").Value = 107.5237638247644
$ws.Range("R5").Value = 967.7138744228799
$ws.Range("S5").Value = 0.002709269186139656
$ws.Range("T5").Value = 0.002709269186139656

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 82.26089466666666
$ws.Range("H6").Value = 246.782684
$ws.Range("I6").Value = 0.2156609604819841
$ws.Range("J6").Value = 0.2156609604819841
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 80.22623699999998
$ws.Range("N6").Value = 240.678711
$ws.Range("O6").Value = 0.77105616682495
$ws.Range("P6").Value = 0.77105616682495
$ws.Range("Q6").Value = 6599.482031360034
$ws.Range("R6").Value = 59395.33828224031
$ws.Range("S6").Value = 0.1662867135230257
$ws.Range("T6").Value = 0.1662867135230257

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 82.26089466666666
$ws.Range("H7").Value = 246.782684
$ws.Range("I7").Value = 0.2156609604819841
$ws.Range("J7").Value = 0.2156609604819841
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 22.51385866666667
$ws.Range("N7").Value = 67.541576
$ws.Range("O7").Value = 0.2163812016255815
$ws.Range("P7").Value = 0.2163812016255815
$ws.Range("Q7").Value = 1852.010156318887
$ws.Range("R7").Value = 16668.09140686999
$ws.Range("S7").Value = 0.04666497777281876
$ws.Range("T7").Value = 0.04666497777281876

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 186.7091623333333
$ws.Range("H8").Value = 560.127487
$ws.Range("I8").Value = 0.489489902130978
$ws.Range("J8").Value = 0.489489902130978
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.307106666666667
$ws.Range("N8").Value = 3.92132
$ws.Range("O8").Value = 0.01256263154946851
$ws.Range("P8").Value = 0.01256263154946851
$ws.Range("Q8").Value = 244.0487908136489
$ws.Range("R8").Value = 2196.43911732284
$ws.Range("S8").Value = 0.006149281287656878
$ws.Range("T8").Value = 0.006149281287656878

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 186.7091623333333
$ws.Range("H9").Value = 560.127487
$ws.Range("I9").Value = 0.489489902130978
$ws.Range("J9").Value = 0.489489902130978
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 80.22623699999998
$ws.Range("N9").Value = 240.678711
$ws.Range("O9").Value = 0.77105616682495
$ws.Range("P9").Value = 0.77105616682495
$ws.Range("Q9").Value = 14978.97350742547
$ws.Range("R9").Value = 134810.7615668292
$ws.Range("S9").Value = 0.3774242076366319
$ws.Range("T9").Value = 0.3774242076366319

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 186.7091623333333
$ws.Range("H10").Value = 560.127487
$ws.Range("I10").Value = 0.489489902130978
$ws.Range("J10").Value = 0.489489902130978
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 22.51385866666667
$ws.Range("N10").Value = 67.541576
$ws.Range("O10").Value = 0.2163812016255815
$ws.Range("P10").Value = 0.2163812016255815
$ws.Range("Q10").Value = 4203.54369254439
$ws.Range("R10").Value = 37831.89323289951
$ws.Range("S10").Value = 0.1059164132066893
$ws.Range("T10").Value = 0.1059164132066893

